$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.368.24'
$ws.Range("E2").Value = '  +2.24%  '
$ws.Range("D3").Value = '2.100.75'
$ws.Range("E3").Value = '  +0.18%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.007'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '344.44'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.34%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.008'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.06%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5224'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +1.58%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4413'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +0.21%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '54.43'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +2.24%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.09337'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +1.58%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.170'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.02%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '24.76'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.67%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '8.648'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +6.04%  '
$ws.Range("D14").Value = '2.107.85'
$ws.Range("E14").Value = '  +0.49%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.906'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +2.28%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '101.58'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +2.36%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001157'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +0.84%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.009'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.15%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '21.15'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +1.95%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.06716'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +1.05%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.375'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +3.16%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.006'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.08%  '
$ws.Range("D23").Value = '30.410.04'
$ws.Range("E23").Value = '  +2.22%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.52'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.42%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.302'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.45%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '21.86'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.22%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '162.51'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.16%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.509'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.26%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '133.22'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +0.55%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.132'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.14%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1050'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.35%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.664'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +1.73%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.725'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +10.82%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.217'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +1.03%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.922'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -1.17%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '10.30'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +1.48%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02631'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +2.84%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06765'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.89%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.7012'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +2.35%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.345'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +4.06%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '12.51'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +1.06%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.2221'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.24%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6820'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +2.66%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.33'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +1.90%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.344'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +2.23%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.007'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.06%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.391'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +19.95%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.648'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +1.17%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00000000354'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +3.93%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.213'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +9.13%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.216'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.07%  '
